# Generate Report for Handback
# ---------------------------------------------------------------
# This script updates the localization-status workbook to reflect
# that the zh-cn and de-de handback packages have been generated:
#   - Status cells move from "Ready for handoff" to
#     "Handed back: in sync with en-US"
#   - The "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated for zh-cn
#     and de-de with the handback file names / timestamps.
#   - A hyperlink is added on the new "Latest Target File" cell,
#     matching the style already used for the source file link.
#   - The now-wider text in those columns is accommodated by
#     widening a few columns.
# ---------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---- Status column: every cell that used to read "Ready for
# ---- handoff" now reads the handed-back status. ----
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$zhcn.Range("C2").Value = $newStatus
$dede.Range("C2").Value = $newStatus

$sourceFileName = "84b984e4-3f7c-4cbb-a06c-fd148656469f.md"
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ede1494eb64b8b36019ec96f40b8058b2ad5c5b8/e2e/84b984e4-3f7c-4cbb-a06c-fd148656469f.md"

# Color/underline used throughout the workbook for hyperlink cells
# (explicit RGB 6495ED, stored as BGR-ish long for the COM Color
# property == RGB(0x64,0x95,0xED)).
$hyperlinkColor = 15570276

# ---- zh-cn sheet row 2 ----
$zhcn.Range("I2").Value = $sourceFileName
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $sourceUrl, "", "", $sourceFileName)
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = $hyperlinkColor

$zhcn.Range("J2").Value = "84b984e4-3f7c-4cbb-a06c-fd148656469f.4d7d596b9d962effe887162d476bbe704ceca03d.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-30 09:27:22"

# ---- de-de sheet row 2 ----
$dede.Range("I2").Value = $sourceFileName
$dede.Hyperlinks.Add($dede.Range("I2"), $sourceUrl, "", "", $sourceFileName)
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = $hyperlinkColor

$dede.Range("J2").Value = "84b984e4-3f7c-4cbb-a06c-fd148656469f.4d7d596b9d962effe887162d476bbe704ceca03d.de-de.xlf"
$dede.Range("K2").Value = "2016-08-30 09:27:29"

# ---- Widen columns to fit the newly-populated, longer text ----
# (Status columns, and the Latest Target File / Latest Handback File
# columns which now hold long .md / .xlf file names.)
$overview.Columns.Item(5).ColumnWidth = 29.17   # Overview!E (zh-cn status)
$overview.Columns.Item(6).ColumnWidth = 29.17   # Overview!F (de-de status)

$zhcn.Columns.Item(3).ColumnWidth = 29.17   # zh-cn!C Status
$zhcn.Columns.Item(9).ColumnWidth = 39.17   # zh-cn!I Latest Target File
$zhcn.Columns.Item(10).ColumnWidth = 39.17  # zh-cn!J Latest Handback File

$dede.Columns.Item(3).ColumnWidth = 29.17   # de-de!C Status
$dede.Columns.Item(9).ColumnWidth = 39.17   # de-de!I Latest Target File
$dede.Columns.Item(10).ColumnWidth = 39.17  # de-de!J Latest Handback File

Write-Host "Handback report generated."
